$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Insert the three new settings rows (keeping the alphabetical ordering
#     used throughout the sheet) -------------------------------------------

# "Charge Energy Density (Wh/kg)" goes right after "Charge Capacity (mWh)"
# and before "CtrCyc" -> final row 8.
$ws.Rows.Item(8).Insert()
$ws.Range("A8").Value = "Charge Energy Density (Wh/kg)"
$ws.Range("B8").Value = "%f"

# "Discharge Energy Density (Wh/kg)" goes right after "Discharge Capacity (mWh)"
# and before "Internal Resistance 1 (mOhm)" -> final row 19.
$ws.Rows.Item(19).Insert()
$ws.Range("A19").Value = "Discharge Energy Density (Wh/kg)"
$ws.Range("B19").Value = "%f"

# "Test" goes right after "Step Time (Seconds)" and before
# "Total Time (Seconds)" -> final row 37.
$ws.Rows.Item(37).Insert()
$ws.Range("A37").Value = "Test"
$ws.Range("B37").Value = "%f"

# --- Colour-code column A (skip the bold header row 1) --------------------

# Default: yellow highlight for every setting row.
$ws.Range("A2:A39").Interior.Color = 65535

# Red highlight for the two brand-new energy-density rows.
$ws.Range("A8").Interior.Color = 255
$ws.Range("A19").Interior.Color = 255

# Green highlight for the counter / PT-100 rows.
$ws.Range("A9:A11").Interior.Color = 5287936
$ws.Range("A29").Interior.Color = 5287936

# --- Refresh the selection / viewport to match the author's saved state ---
$ws.Range("A14").Select()
